# Daily attendance processing - reorder "Recorded By" (column G) entries.
#
# Each G cell holds a comma-separated list of recorders, e.g.
#   "System, dnasr281@gmail.com"  or  "system, backup@backdoor.com, System"
#
# Rule observed in the source data: the literal, case-sensitive token
# "System" (if it is the LAST entry) stays pinned at the end of the list;
# every other entry in the cell has its order reversed. If "System" is not
# the last entry (or the cell has only one entry), the whole list is simply
# reversed.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count
$colG = 7

for ($r = 1; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, $colG)
    $raw = $cell.Value2

    if ($null -eq $raw) {
        continue
    }
    if ($raw -eq "") {
        continue
    }

    $parts = $raw -split ", "

    if ($parts.Length -le 1) {
        continue
    }

    $lastIdx = $parts.Length - 1

    if ($parts[$lastIdx].Equals("System")) {
        if ($parts.Length -eq 2) {
            $newParts = @($parts[0], "System")
        } else {
            $head = $parts[0..($lastIdx - 1)]
            $headRev = $head[($head.Length - 1)..0]
            $newParts = $headRev + @("System")
        }
    } else {
        $newParts = $parts[$lastIdx..0]
    }

    $newVal = [string]::Join(", ", $newParts)

    if (-not $newVal.Equals($raw)) {
        $cell.Value = $newVal
    }
}
